# Updates the cryptos list price/volume figures in Sheet1 (columns D & E).
# Values that look like plain numbers are prefixed with a leading
# apostrophe so Excel keeps them as text (matching the source sheet's
# existing inline-string/text formatting for the Price column).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.025.57"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "3.517.01"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'592.80"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("D6").Value = "'133.88"
$ws.Range("E6").Value = "  -1.46%  "
$ws.Range("D7").Value = "3.516.35"
$ws.Range("E7").Value = "  -1.47%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.489"
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("D10").Value = "'0.124"
$ws.Range("E10").Value = "  +0.86%  "
$ws.Range("D11").Value = "'7.18"
$ws.Range("E11").Value = "  +2.79%  "
$ws.Range("D12").Value = "'0.385"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "4.118.93"
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("D14").Value = "'27.70"
$ws.Range("E14").Value = "  +2.08%  "
$ws.Range("D15").Value = "'0.0000181"
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").Value = "3.516.61"
$ws.Range("E17").Value = "  -1.39%  "
$ws.Range("D18").Value = "65.012.49"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").Value = "'10.04"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "'14.31"
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").Value = "'5.67"
$ws.Range("E21").Value = "  -3.32%  "
$ws.Range("D22").Value = "'391.86"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("D23").Value = "'0.576"
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("D24").Value = "3.659.02"
$ws.Range("E24").Value = "  -1.46%  "
$ws.Range("D25").Value = "'74.63"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  -4.04%  "
$ws.Range("E28").Value = "  +8.59%  "
$ws.Range("D29").Value = "'7.65"
$ws.Range("E29").Value = "  -1.79%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").Value = "'2.26"
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("D32").Value = "'8.29"
$ws.Range("E32").Value = "  -2.16%  "
$ws.Range("D33").Value = "3.523.28"
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("D34").Value = "'24.07"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").Value = "'5.24"
$ws.Range("E37").Value = "  +4.32%  "
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("D39").Value = "'6.93"
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("D40").Value = "'168.97"
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("D41").Value = "'0.0805"
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("E42").Value = "  -0.78%  "
$ws.Range("D43").Value = "'1.25"
$ws.Range("E43").Value = "  +3.45%  "
$ws.Range("D44").Value = "'25.85"
$ws.Range("E44").Value = "  -5.53%  "
$ws.Range("D45").Value = "'42.91"
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").Value = "'4.43"
$ws.Range("E47").Value = "  -1.08%  "
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("D49").Value = "'6.88"
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("D50").Value = "2.432.71"
$ws.Range("E50").Value = "  -2.25%  "
$ws.Range("D51").Value = "'0.901"
$ws.Range("E51").Value = "  +3.03%  "
